# Generate Report for Handback
# Update the "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# timestamps on the per-locale report sheets to reflect the latest handback run.

$wb = $excel.ActiveWorkbook

# zh-cn sheet: rows 4 and 5 share the same handoff/handback timestamps
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E4").Value = "2016-03-23 12:23:23"
$wsZhCn.Range("E5").Value = "2016-03-23 12:23:23"
$wsZhCn.Range("H4").Value = "2016-03-23 12:23:50"
$wsZhCn.Range("H5").Value = "2016-03-23 12:23:50"

# de-de sheet: rows 4 and 5 share the same handoff/handback timestamps
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E4").Value = "2016-03-23 12:23:27"
$wsDeDe.Range("E5").Value = "2016-03-23 12:23:27"
$wsDeDe.Range("H4").Value = "2016-03-23 12:23:57"
$wsDeDe.Range("H5").Value = "2016-03-23 12:23:57"
